$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were deleted upstream (RM 232, SC 92),
# which shifts all subsequent rows up by one each.
$ws.Rows("26:26").Delete()  # was "RM 232"
$ws.Rows("27:27").Delete()  # was "SC 92" (now at row 27 after first delete)

# Apply the remaining per-cell value changes (post row-shift).
$ws.Range("C2").Value = 14.9
$ws.Range("F2").Value = 18.03
$ws.Range("E3").ClearContents()
$ws.Range("E4").Value = -6.4
$ws.Range("C6").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("F9").ClearContents()
$ws.Range("C12").Value = 12.5
$ws.Range("F13").Value = 17.1
$ws.Range("C14").ClearContents()
$ws.Range("E15").Value = -8.4
$ws.Range("E18").Value = -8.5
$ws.Range("E19").ClearContents()
$ws.Range("F19").ClearContents()
$ws.Range("C20").Value = 12.5
$ws.Range("C21").Value = 12.7
$ws.Range("E22").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("E23").Value = -7
$ws.Range("C24").ClearContents()
$ws.Range("E25").Value = -7.1
$ws.Range("F25").Value = 16.6
$ws.Range("B26").Value = -20.2
$ws.Range("B27").ClearContents()
$ws.Range("E27").ClearContents()
$ws.Range("B28").ClearContents()
$ws.Range("F28").Value = 17.44
$ws.Range("B29").Value = -19.5
$ws.Range("B30").Value = -19.7
$ws.Range("B31").ClearContents()
$ws.Range("C31").Value = 15.3
$ws.Range("F31").ClearContents()
$ws.Range("B32").ClearContents()
$ws.Range("F32").Value = 17.39
$ws.Range("C33").Value = 10.4
